$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.84"
$ws.Range("E2").Value = "'0.20%"
$ws.Range("D3").Value = "'36.81"
$ws.Range("E3").Value = "'3.46%"
$ws.Range("D4").Value = "'5.026"
$ws.Range("E4").Value = "'-1.26%"
$ws.Range("D5").Value = "'0.07834"
$ws.Range("E5").Value = "'0.10%"
$ws.Range("D6").Value = "'2.170"
$ws.Range("E6").Value = "'-3.87%"
$ws.Range("D7").Value = "'8.026"
$ws.Range("E7").Value = "'-1.16%"
$ws.Range("D8").Value = "'4.057"
$ws.Range("E8").Value = "'1.29%"
$ws.Range("D9").Value = "'0.9239"
$ws.Range("E9").Value = "'-0.12%"
$ws.Range("D10").Value = "'0.09953"
$ws.Range("E10").Value = "'2.61%"
$ws.Range("D11").Value = "'0.1869"
$ws.Range("E11").Value = "'2.48%"
$ws.Range("D12").Value = "'0.08714"
$ws.Range("E12").Value = "'-0.21%"
$ws.Range("D13").Value = "'0.03587"
$ws.Range("E13").Value = "'4.98%"
$ws.Range("D14").Value = "'0.09941"
$ws.Range("E14").Value = "'-0.03%"
$ws.Range("D15").Value = "'0.001480"
$ws.Range("E15").Value = "'-0.23%"
$ws.Range("D16").Value = "'0.005670"
$ws.Range("E16").Value = "'-0.14%"
$ws.Range("D17").Value = "'3.465"
$ws.Range("E17").Value = "'-0.55%"
$ws.Range("D18").Value = "'2.340"
$ws.Range("E18").Value = "'8.92%"
$ws.Range("D19").Value = "'0.3447"
$ws.Range("E19").Value = "'0.80%"
$ws.Range("E20").Value = "'1.89%"
$ws.Range("D21").Value = "'4.916"
$ws.Range("E21").Value = "'8.46%"
$ws.Range("E22").Value = "'-1.42%"
$ws.Range("D23").Value = "'0.04599"
$ws.Range("E23").Value = "'-1.67%"
$ws.Range("D24").Value = "'0.005177"
$ws.Range("E24").Value = "'13.85%"
$ws.Range("D25").Value = "'0.001233"
$ws.Range("E25").Value = "'-0.73%"
$ws.Range("E26").Value = "'7.93%"
$ws.Range("D27").Value = "'0.0002719"
$ws.Range("E27").Value = "'0.83%"
$ws.Range("D39").Value = "'0.01814"
$ws.Range("E39").Value = "'3.21%"
$ws.Range("D40").Value = "'0.04736"
$ws.Range("E40").Value = "'0.55%"
$ws.Range("D41").Value = "'0.007920"
$ws.Range("E41").Value = "'-0.51%"
$ws.Range("E42").Value = "'-1.06%"
$ws.Range("D43").Value = "'0.007595"
$ws.Range("E43").Value = "'-5.13%"
$ws.Range("D44").Value = "'0.002232"
$ws.Range("E44").Value = "'-2.83%"
$ws.Range("D45").Value = "'0.01046"
$ws.Range("E45").Value = "'14.70%"
$ws.Range("D46").Value = "'0.00006333"
$ws.Range("E46").Value = "'1.82%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.22%"
$ws.Range("D48").Value = "'0.0005803"
$ws.Range("E48").Value = "'0.04%"
$ws.Range("D49").Value = "'33.50"
$ws.Range("E49").Value = "'727.25%"
$ws.Range("E50").Value = "'0.22%"
$ws.Range("D51").Value = "'0.00002102"
$ws.Range("E51").Value = "'0.22%"
